$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column D ---
# Existing data in columns D:K (prior fiscal years) shifts right to E:L,
# making room for a new "latest fiscal year" column at D.
$ws.Columns("D:D").Insert()

# --- Copy formatting (number format / bold / alignment) from column E into the new column D ---
for ($r = 7; $r -le 102; $r++) {
  $srcCell = $ws.Cells.Item($r, 5)
  $dstCell = $ws.Cells.Item($r, 4)
  $dstCell.NumberFormat = $srcCell.NumberFormat
  $dstCell.Font.Bold = $srcCell.Font.Bold
  $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
}

# --- Adjust column widths to account for the extra data column ---
$ws.Range("A:A").ColumnWidth = 5.78
$ws.Range("B:B").ColumnWidth = 26.01
$ws.Range("C:C").ColumnWidth = 68.23
$ws.Range("D:K").ColumnWidth = 13.78
$ws.Range("L:L").ColumnWidth = 8.23

# --- Populate the new column D with the latest (most recent) fiscal-period figures ---
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 718900
$ws.Cells.Item(9, 4).Value = 353300
$ws.Cells.Item(10, 4).Value = 365600
$ws.Cells.Item(12, 4).Value = 76000
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 4200
$ws.Cells.Item(15, 4).Value = 13600
$ws.Cells.Item(17, 4).Value = 547300
$ws.Cells.Item(18, 4).Value = 171600
$ws.Cells.Item(20, 4).Value = 800
$ws.Cells.Item(21, 4).Value = 186000
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 4).Value = 172400
$ws.Cells.Item(24, 4).Value = 19500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 152900
$ws.Cells.Item(27, 4).Value = 152800
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = -5700
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -800
$ws.Cells.Item(33, 4).Value = 147000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 147000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 349300
$ws.Cells.Item(42, 4).Value = 2500
$ws.Cells.Item(43, 4).Value = 102700
$ws.Cells.Item(44, 4).Value = 98000
$ws.Cells.Item(45, 4).Value = 16000
$ws.Cells.Item(46, 4).Value = 568400
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 31300
$ws.Cells.Item(49, 4).Value = 156800
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 60000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 816500
$ws.Cells.Item(57, 4).Value = 39600
$ws.Cells.Item(58, 4).Value = "NA"
$ws.Cells.Item(59, 4).Value = 70700
$ws.Cells.Item(60, 4).Value = 110300
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(62, 4).Value = 88200
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 209700
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 512800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 606800
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 147000
$ws.Cells.Item(83, 4).Value = 13600
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 146000
$ws.Cells.Item(91, 4).Value = -20300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -113600
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -97100
$ws.Cells.Item(101, 4).Value = -1000
$ws.Cells.Item(102, 4).Value = -65700
